# Updated symbol list on Fri Dec 23 13:25:58 UTC 2022 with GitHub Actions
#
# Applies the "One" coin re-ranking (moves from row 18 up to row 10, pushing
# WazirX..CoinExToken down by one row) plus the various Price (column D)
# refreshes and the single Volume(1h) text tweak on row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $rng = $ws.Range($Addr)
    # Force text storage so numeric-looking strings (e.g. "246.15") are not
    # silently reinterpreted as floating point numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $Val
}

# ---- Simple Price (column D) refreshes, rows 2-9 ----
Set-TextValue "D2" "246.15"
Set-TextValue "D3" "22.03"
Set-TextValue "D4" "5.419"
Set-TextValue "D5" "0.05869"
Set-TextValue "D6" "3.389"
Set-TextValue "D7" "6.359"
Set-TextValue "D8" "0.8142"
Set-TextValue "D9" "1.029"

# ---- Rows 10-18: "One" jumps from row 18 to row 10; WazirX..CoinExToken ----
# ---- each shift down by one row, with refreshed price/volume data.     ----

# Row 10 (was WazirX, now One)
Set-TextValue "B10" "One"
Set-TextValue "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.0005889"
Set-TextValue "E10" "9OneONE"

# Row 11 (was MandalaExchangeToken, now WazirX)
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1428"
Set-TextValue "E11" "10WazirXWRX"

# Row 12 (was LiechtensteinCryptoassetsExchange, now MandalaExchangeToken)
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07452"
Set-TextValue "E12" "11MandalaExchangeTokenMDX"

# Row 13 (was BitrueCoin, now LiechtensteinCryptoassetsExchange)
Set-TextValue "B13" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03466"
Set-TextValue "E13" "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14 (was MCDex, now BitrueCoin)
Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03022"
Set-TextValue "E14" "13BitrueCoinBTR"

# Row 15 (was BitMartToken, now MCDex)
Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "4.207"
Set-TextValue "E15" "14MCDexMCB"

# Row 16 (was BitForexToken, now BitMartToken)
Set-TextValue "B16" "BitMartToken"
Set-TextValue "C16" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09398"
Set-TextValue "E16" "15BitMartTokenBMX"

# Row 17 (was CoinExToken, now BitForexToken)
Set-TextValue "B17" "BitForexToken"
Set-TextValue "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001595"
Set-TextValue "E17" "16BitForexTokenBF"

# Row 18 (was One, now CoinExToken)
Set-TextValue "B18" "CoinExToken"
Set-TextValue "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04824"
Set-TextValue "E18" "17CoinExTokenCET"

# ---- Remaining Price (column D) refreshes, rows 19-25 ----
Set-TextValue "D19" "0.006086"
Set-TextValue "D20" "0.004123"
Set-TextValue "D21" "0.0009980"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.696"
Set-TextValue "D24" "2.220"
Set-TextValue "D25" "0.3245"

# ---- Row 27: Volume(1h) text tweak (no price change) ----
Set-TextValue "E27" "26UpBotsUBXTWorstin24h"

# ---- Remaining Price (column D) refreshes, rows 40-49 ----
Set-TextValue "D40" "0.03858"
Set-TextValue "D41" "0.006638"
Set-TextValue "D43" "0.002599"
Set-TextValue "D44" "0.006690"
Set-TextValue "D45" "0.00005621"
Set-TextValue "D48" "0.1408"
Set-TextValue "D49" "0.00002100"
